$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format for Price cells whose new values would otherwise be
# auto-converted to numbers by Excel (losing exact text representation).
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D51').NumberFormat = '@'

# Apply updated cell values row by row
$ws.Range('D2').Value = '26.241.31'
$ws.Range('E2').Value = '  -1.37%  '
$ws.Range('D3').Value = '1.672.50'
$ws.Range('E3').Value = '  +1.14%  '
$ws.Range('D4').Value = '1.006'
$ws.Range('E4').Value = '  -0.27%  '
$ws.Range('D5').Value = '217.60'
$ws.Range('E5').Value = '  -1.07%  '
$ws.Range('D6').Value = '0.5136'
$ws.Range('E6').Value = '  +0.89%  '
$ws.Range('E7').Value = '  -0.27%  '
$ws.Range('D8').Value = '0.2659'
$ws.Range('E8').Value = '  +5.16%  '
$ws.Range('D9').Value = '0.06383'
$ws.Range('E9').Value = '  +4.33%  '
$ws.Range('D10').Value = '21.57'
$ws.Range('E10').Value = '  -0.18%  '
$ws.Range('D11').Value = '0.07390'
$ws.Range('E11').Value = '  +0.59%  '
$ws.Range('D12').Value = '1.683.22'
$ws.Range('E12').Value = '  +1.73%  '
$ws.Range('D13').Value = '4.551'
$ws.Range('E13').Value = '  +2.62%  '
$ws.Range('D14').Value = '0.5830'
$ws.Range('E14').Value = '  +1.79%  '
$ws.Range('D15').Value = '1.901.16'
$ws.Range('E15').Value = '  +1.25%  '
$ws.Range('D16').Value = '0.000008701'
$ws.Range('E16').Value = '  +8.25%  '
$ws.Range('D17').Value = '64.74'
$ws.Range('E17').Value = '  +0.40%  '
$ws.Range('D18').Value = '26.303.89'
$ws.Range('E18').Value = '  -1.23%  '
$ws.Range('D19').Value = '4.956'
$ws.Range('E19').Value = '  +0.14%  '
$ws.Range('D20').Value = '1.006'
$ws.Range('E20').Value = '  -0.28%  '
$ws.Range('E21').Value = '  +2.98%  '
$ws.Range('D22').Value = '189.36'
$ws.Range('E22').Value = '  +4.62%  '
$ws.Range('D23').Value = '6.220'
$ws.Range('E23').Value = '  +0.54%  '
$ws.Range('D24').Value = '1.007'
$ws.Range('E24').Value = '  -0.44%  '
$ws.Range('D25').Value = '144.57'
$ws.Range('E25').Value = '  +1.13%  '
$ws.Range('D26').Value = '7.630'
$ws.Range('E26').Value = '  +1.01%  '
$ws.Range('D27').Value = '0.1186'
$ws.Range('E27').Value = '  +3.42%  '
$ws.Range('D28').Value = '15.63'
$ws.Range('E28').Value = '  +4.22%  '
$ws.Range('D29').Value = '0.05964'
$ws.Range('E29').Value = '  +2.74%  '
$ws.Range('D30').Value = '1.283'
$ws.Range('E30').Value = '  -3.77%  '
$ws.Range('D31').Value = '1.319'
$ws.Range('E31').Value = '  -1.48%  '
$ws.Range('D32').Value = '3.528'
$ws.Range('E32').Value = '  +3.56%  '
$ws.Range('D33').Value = '3.528'
$ws.Range('E33').Value = '  +3.61%  '
$ws.Range('E34').Value = '  +4.20%  '
$ws.Range('D35').Value = '1.015'
$ws.Range('E35').Value = '  +4.18%  '
$ws.Range('D36').Value = '0.6027'
$ws.Range('E36').Value = '  +1.71%  '
$ws.Range('E37').Value = '  -2.38%  '
$ws.Range('D39').Value = '0.01620'
$ws.Range('E39').Value = '  +2.87%  '
$ws.Range('D40').Value = '6.089'
$ws.Range('E40').Value = '  +6.01%  '
$ws.Range('D41').Value = '1.079.12'
$ws.Range('E41').Value = '  +1.02%  '
$ws.Range('D42').Value = '0.8703'
$ws.Range('E42').Value = '  +0.81%  '
$ws.Range('D43').Value = '1.010'
$ws.Range('E43').Value = '  -0.10%  '
$ws.Range('D44').Value = '100.11'
$ws.Range('E44').Value = '  +4.65%  '
$ws.Range('D45').Value = '1.821.68'
$ws.Range('E45').Value = '  +1.65%  '
$ws.Range('D46').Value = '0.00000000112'
$ws.Range('E46').Value = '  +6.76%  '
$ws.Range('D47').Value = '56.12'
$ws.Range('E47').Value = '  +1.80%  '
$ws.Range('E48').Value = '  +0.47%  '
$ws.Range('D49').Value = '8.053'
$ws.Range('E49').Value = '  +4.01%  '
$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D50').Value = '0.05208'
$ws.Range('E50').Value = '  +0.04%  '
$ws.Range('B51').Value = 'Mantle'
$ws.Range('C51').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D51').Value = '0.4300'
$ws.Range('E51').Value = '  -1.80%  '
